$d = $word.ActiveDocument

function Merge-InParagraph {
    param($ParaIndex, $OldText, $NewText, $ForceColor)
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $p.Range
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if ($ForceColor -ne $null) {
        $p2 = $d.Paragraphs.Item($ParaIndex)
        $r2 = $p2.Range
        $r2.Find.ClearFormatting()
        $found = $r2.Find.Execute($NewText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $r2.Font.Color = $ForceColor
        }
    }
    return $ok
}

# --- 1. TOC hyperlink "2 Hi" + "gh Concept" -> "2 High Concept" (paragraph 63) ---
$old1 = "2 Hi" + "gh Concept"
$new1 = "2 High Concept"
Merge-InParagraph 63 $old1 $new1 $null | Out-Null

# --- 2. TOC hyperlink "11.1 Ch" + "aracters" -> "11.1 Characters" (paragraph 77) ---
# Force the text color back to the hyperlink blue (1155CC -> wdColor 0xCC5511) since the
# paragraph also contains a leading "  " run colored 222222 that Find/Replace formatting
# could otherwise bleed into the merged run.
$old2 = "11.1 Ch" + "aracters"
$new2 = "11.1 Characters"
$blue = 13391121
Merge-InParagraph 77 $old2 $new2 $blue | Out-Null

# --- 3. "Rat" + "ing: " -> "Rating: " (paragraph 148) ---
$old3 = "Rat" + "ing: "
$new3 = "Rating: "
Merge-InParagraph 148 $old3 $new3 $null | Out-Null

# --- 4. "Publisher: " + "brandon's" + " games" -> "Publisher: GBB (games by brandon)" (paragraph 151) ---
$apos = [char]8217
$old4 = "brandon" + $apos + "s games"
$new4 = "GBB (games by brandon)"
Merge-InParagraph 151 $old4 $new4 $null | Out-Null

# --- 5. Heading "2 Hig" + "h Concept" -> "2 High Concept" (paragraph 154) ---
$old5 = "2 Hig" + "h Concept"
$new5 = "2 High Concept"
Merge-InParagraph 154 $old5 $new5 $null | Out-Null

# --- 6. Heading "5" + " Competitors / Similar Titles" -> "5 Competitors / Similar Titles" (paragraph 162) ---
$old6 = "5" + " Competitors / Similar Titles"
$new6 = "5 Competitors / Similar Titles"
Merge-InParagraph 162 $old6 $new6 $null | Out-Null

# --- 7. Heading "1" + "0 Game Play" -> "10 Game Play" (paragraph 173) ---
$old7 = "1" + "0 Game Play"
$new7 = "10 Game Play"
Merge-InParagraph 173 $old7 $new7 $null | Out-Null

# --- 8. "1" + "0.2 " -> "10.2 " (paragraph 176) ---
$old8 = "1" + "0.2 "
$new8 = "10.2 "
Merge-InParagraph 176 $old8 $new8 $null | Out-Null

# --- 9. Heading "1" + "1 Players" -> "11 Players" (paragraph 178) ---
$old9 = "1" + "1 Players"
$new9 = "11 Players"
Merge-InParagraph 178 $old9 $new9 $null | Out-Null

# --- 10. "Character" + "s" -> "Characters" (paragraph 179) ---
$old10 = "Character" + "s"
$new10 = "Characters"
Merge-InParagraph 179 $old10 $new10 $null | Out-Null

# --- 11. "1" + "3.1 " -> "13.1 " (paragraph 189) ---
$old11 = "1" + "3.1 "
$new11 = "13.1 "
Merge-InParagraph 189 $old11 $new11 $null | Out-Null

# --- 12. Add (empty) even/default/first-page footers, matching the header setup ---
$sec = $d.Sections(1)
$footers = $sec.Footers
$footers.Item(1).Range.Text = ""
$footers.Item(2).Range.Text = ""
$footers.Item(3).Range.Text = ""

Write-Output "done"
